$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-22 10:28:16"
$zhcn.Range("K2").Value = "2016-08-22 10:28:33"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-22 10:28:21"
$dede.Range("K2").Value = "2016-08-22 10:28:39"
